# Correção nos dados e inicio da analise PNAD 2009
#
# The worksheet "dados" contains a table where two section-header rows
# ("situação do domicílio" at row 5 and "grandes regiões e unidades da
# federação" at row 8) have no numeric data. Deleting those two rows shifts
# every subsequent row up, which both realigns the labels with their correct
# data and removes the two now-unused trailing rows (39 and 40).
#
# Also, the column-2 header label (B2) is renamed from the placeholder
# "unnamed: 1_level_1" to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the placeholder header text.
$ws.Range("B2").Value = "total"

# Delete row 8 first (higher row number) so row 5's index doesn't shift
# before we get to delete it.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
